$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "time" column (E2:E5) from 15:44 to 16:16, keeping the value as plain text
$ws.Range("E2").Value = "16:16"
$ws.Range("E3").Value = "16:16"
$ws.Range("E4").Value = "16:16"
$ws.Range("E5").Value = "16:16"
